$d = $word.ActiveDocument
$d.Content.Find.Execute("NOM  ", $true, $false, $false, $false, $false, $true, 1, $false, "nom  ", 2)
